$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.261.70"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.307.96"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.15"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "558.13"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "3.300.47"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.589"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.10"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.71"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.834.58"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "615.98"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.10"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "66.222.87"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "3.311.44"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.11"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.913"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.31"
$ws.Range("E23").Value = "  +7.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.00"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.98"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.95"
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.00"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.77"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.79"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.71"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.36"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.11"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.75"
$ws.Range("E33").Value = "  +6.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "564.73"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.12"
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "3.776.79"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.54"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.48"
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("D41").Value = "0.0₃0728"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.130"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.340"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0426"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.03%  "
